$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("in")

# New rows to append (A: id, B: slug, C: name)
$newRows = @(
    @("ID:ZN0415", "valentis", "Valentis"),
    @("ID:ZN0414", "la-vida-caribena", "La Vida Caribena"),
    @("ID:ZN0411", "renu", "ReNu"),
    @("ID:ZN0412", "ecover", "Ecover"),
    @("ID:13928", "relax", "Relax"),
    @("ID:ZN0410", "lip-smacker", "Lip Smacker"),
    @("ID:ZN0408", "bloom-robbins", "Bloom Robbins"),
    @("ID:ZN0413", "dr-althea", "Dr. Althea"),
    @("ID:ZN0409", "jovo", "JöVö"),
    @("ID:43585", "satin-care", "Satin Care"),
    @("ID:ZN0416", "gallus", "Gallus"),
    @("ID:ZN0417", "fackelmann", "Fackelmann"),
    @("ID:ZN0418", "fammilky", "Fammilky"),
    @("ID:ZN0419", "flawless", "Flawless"),
    @("ID:ZN0420", "q-svice-ozona", "Q svíce Ozóna")
)

$startRow = 813
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Adjust column widths to fit new content (matches bestFit behaviour)
$ws.Range("A1:C827").EntireColumn.AutoFit() | Out-Null

# Reset selection / view to top-left cell, no special scroll/selection
$ws.Range("A1").Select()
